# Applies the edit described by the commit diff:
#  - Adds a new "STEVE-O" transformer row (crudeoil -> gasoline 60% + kerosene 40%)
#    to the Transformers sheet, and splits the existing Refinery row's output
#    into gasoline/jetfuel/kerosene fractions.
#  - Makes the Transformers sheet the active tab (was Connectors).

$wb = $excel.ActiveWorkbook

$wsTransformers = $wb.Worksheets.Item("Transformers")

# --- Update existing row 2 (Refinery) on the Transformers sheet ---
# G2 goes from 1 -> 0.5, and new product/efficiency pairs are appended.
$wsTransformers.Range("G2").Value = 0.5
$wsTransformers.Range("H2").Value = "jetfuel"
$wsTransformers.Range("I2").Value = 0.3
$wsTransformers.Range("J2").Value = "kerosene"
$wsTransformers.Range("K2").Value = 0.2

# --- Add a brand-new row 4 (STEVE-O) on the Transformers sheet ---
$wsTransformers.Range("A4").Value = "STEVE-O"
$wsTransformers.Range("B4").Value = "crudeoil"
$wsTransformers.Range("C4").Value = 0
$wsTransformers.Range("D4").Value = 0
$wsTransformers.Range("E4").Value = 0.95
$wsTransformers.Range("F4").Value = "gasoline"
$wsTransformers.Range("G4").Value = 0.6
$wsTransformers.Range("H4").Value = "kerosene"
$wsTransformers.Range("I4").Value = 0.4

# --- Column width tweaks on the Transformers sheet to fit the new data ---
# (target stored widths ~6.43 / ~5.71 / ~9.29 / ~9.29 chars; the values below
#  are the closest settable ColumnWidth inputs for those targets)
$wsTransformers.Columns.Item(3).ColumnWidth = 5.6666666666667
$wsTransformers.Columns.Item(4).ColumnWidth = 4.8333333333333
$wsTransformers.Columns.Item(8).ColumnWidth = 8.5
$wsTransformers.Columns.Item(10).ColumnWidth = 8.5

# --- Make Transformers the active sheet (was Connectors) ---
$wsTransformers.Activate()

# --- Selection on Transformers moves to R14 ---
$wsTransformers.Range("R14").Select() | Out-Null
